$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of Price (column D) values are plain decimal numbers (e.g.
# "598.20", "10.00", "0.0000290"). If written into a General-formatted
# cell, Excel auto-converts numeric-looking text into a real number and
# loses the exact text (trailing zeros, tiny decimals, etc). The source
# keeps these as literal text, so force a Text number format on just
# those cells before writing their new values.
$textCells = @(
    "D5", "D6", "D11", "D12", "D13", "D14", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D36", "D37", "D39", "D41", "D46", "D51"
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '72.281.78'
$ws.Range('E2').Value = '  +4.40%  '
$ws.Range('D3').Value = '3.623.41'
$ws.Range('E3').Value = '  +7.00%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '598.20'
$ws.Range('E5').Value = '  +1.94%  '
$ws.Range('D6').Value = '182.62'
$ws.Range('E6').Value = '  +1.42%  '
$ws.Range('D7').Value = '3.612.89'
$ws.Range('E7').Value = '  +6.97%  '
$ws.Range('E8').Value = '  +2.03%  '
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('E10').Value = '  +5.20%  '
$ws.Range('D11').Value = '0.610'
$ws.Range('E11').Value = '  +2.94%  '
$ws.Range('D12').Value = '50.72'
$ws.Range('E12').Value = '  +4.39%  '
$ws.Range('D13').Value = '0.0000290'
$ws.Range('E13').Value = '  +2.79%  '
$ws.Range('D14').Value = '705.90'
$ws.Range('E14').Value = '  +3.76%  '
$ws.Range('D15').Value = '4.204.39'
$ws.Range('E15').Value = '  +7.13%  '
$ws.Range('E16').Value = '  +3.98%  '
$ws.Range('D17').Value = '72.399.88'
$ws.Range('E17').Value = '  +4.47%  '
$ws.Range('D18').Value = '3.610.15'
$ws.Range('E18').Value = '  +6.06%  '
$ws.Range('E19').Value = '  +2.03%  '
$ws.Range('D20').Value = '18.63'
$ws.Range('E20').Value = '  +5.05%  '
$ws.Range('D21').Value = '11.74'
$ws.Range('E21').Value = '  +3.81%  '
$ws.Range('D22').Value = '0.934'
$ws.Range('E22').Value = '  +3.34%  '
$ws.Range('E23').Value = '  +8.67%  '
$ws.Range('D24').Value = '17.86'
$ws.Range('E24').Value = '  +4.26%  '
$ws.Range('D25').Value = '105.51'
$ws.Range('E25').Value = '  +2.55%  '
$ws.Range('D26').Value = '4.04'
$ws.Range('E26').Value = '  +2.80%  '
$ws.Range('D27').Value = '2.86'
$ws.Range('E27').Value = '  +4.89%  '
$ws.Range('D28').Value = '10.00'
$ws.Range('E28').Value = '  +4.10%  '
$ws.Range('D29').Value = '35.59'
$ws.Range('E29').Value = '  +4.99%  '
$ws.Range('D30').Value = '9.16'
$ws.Range('E30').Value = '  +4.59%  '
$ws.Range('D31').Value = '7.44'
$ws.Range('E31').Value = '  +6.86%  '
$ws.Range('D32').Value = '4.17'
$ws.Range('E32').Value = '  +13.78%  '
$ws.Range('D33').Value = '595.24'
$ws.Range('E33').Value = '  +7.22%  '
$ws.Range('E34').Value = '  +2.11%  '
$ws.Range('E35').Value = '  +1.73%  '
$ws.Range('D36').Value = '59.65'
$ws.Range('E36').Value = '  +1.90%  '
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '3.646.83'
$ws.Range('E38').Value = '  -0.79%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = '0.145'
$ws.Range('E39').Value = '  +3.34%  '
$ws.Range('D40').Value = '0.0₃0779'
$ws.Range('E40').Value = '  +9.01%  '
$ws.Range('D41').Value = '35.95'
$ws.Range('E41').Value = '  +0.34%  '
$ws.Range('E42').Value = '  +7.04%  '
$ws.Range('E43').Value = '  +4.48%  '
$ws.Range('E44').Value = '  +6.47%  '
$ws.Range('E45').Value = '  +2.95%  '
$ws.Range('D46').Value = '3.44'
$ws.Range('E46').Value = '  +3.67%  '
$ws.Range('E47').Value = '  +4.65%  '
$ws.Range('E48').Value = '  +5.33%  '
$ws.Range('E49').Value = '  +1.83%  '
$ws.Range('E50').Value = '  -0.27%  '
$ws.Range('D51').Value = '134.03'
$ws.Range('E51').Value = '  +0.60%  '
